# Module 1 and Module 2 Question Table Template Extract
#
# 1) Rename the existing sheet and add a second sheet holding a learner
#    user-data table.
# 2) Update several question rows / answer values on the first sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 cell value updates -------------------------------------------

# Row 2 - 6 Sigma price question
$ws1.Range("D2").Value = 20
$ws1.Range("E2").Value = 25
$ws1.Range("F2").Value = 15
$ws1.Range("G2").Value = 5

# Row 3 - internal stakeholders question
$ws1.Range("E3").Value = 4
$ws1.Range("F3").Value = 8
$ws1.Range("G3").Value = 6

# Row 4 - process desks question
$ws1.Range("B4").Value = "How many process desks are a part of your Scope ?"
$ws1.Range("D4").Value = 5
$ws1.Range("E4").Value = 4
$ws1.Range("F4").Value = 8
$ws1.Range("G4").Value = 6

# Row 8 - penalty question (the "No Penalty" text is introduced here,
# ahead of row 5's text, to match the author's original edit order)
$ws1.Range("D8").Value = "No Penalty"
$ws1.Range("F8").Value = 10000
$ws1.Range("G8").Value = 20000

# Row 5 - replaced question (service levels penalty)
$ws1.Range("B5").Value = "Which among the service levels mentioned does not attract Service penalty ?"
$ws1.Range("C5").Value = "%"
$ws1.Range("D5").Value = 95
$ws1.Range("E5").Value = 90
$ws1.Range("F5").Value = 95.01
$ws1.Range("G5").Value = 85
$ws1.Range("H5").Value = 95.01
$ws1.Rows.Item(5).RowHeight = 45

# Row 6 - 4 Sigma price question
$ws1.Range("D6").Value = 30
$ws1.Range("E6").Value = 25

# Row 7 - external stakeholders question
$ws1.Range("D7").Value = 3
$ws1.Range("E7").Value = 1
$ws1.Range("F7").Value = 4
$ws1.Range("G7").Value = 8

# Row 9 - bonus reward question
$ws1.Range("E9").Value = 25
$ws1.Range("F9").Value = 10
$ws1.Range("G9").Value = 20

# --- Sheet names -----------------------------------------------------------

$ws1.Name = "OEModule1Question"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LearnerUserDataTable"

# --- Sheet2 content: learner user-data tracking table ----------------------
# Values are written in the same order the original author entered them, so
# that the regenerated shared-strings table lines up index-for-index.

$ws2.Range("A3").Value = "TableName"

$ws2.Range("A4").Value = "Attempt"
$ws2.Range("B4").Value = "UserEmail"
$ws2.Range("C4").Value = "ProductName"
$ws2.Range("D4").Value = "ModuleName"

$ws2.Range("B3").Value = "OEModuleDiscoverLearnerUserData"

$ws2.Range("E4").Value = "UI-1"
$ws2.Range("F4").Value = "UI-2"
$ws2.Range("G4").Value = "UI-3"
$ws2.Range("H4").Value = "UI-4"

$ws2.Range("M4").Value = "Discover Score"

$ws2.Range("I4").Value = "UI-5"
$ws2.Range("J4").Value = "UI-6"
$ws2.Range("K4").Value = "UI-7"
$ws2.Range("L4").Value = "UI-8"

# Apply the bold "Arial Black" header formatting (matches Sheet1's header
# row style) to row 4, then the same style in red to the A3:B3 table-name
# banner.
$ws1.Range("A1").Copy()
$ws2.Range("A4:M4").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A3:B3").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A3:B3").Font.Color = 255

$ws2.Rows.Item(3).RowHeight = 18.75
$ws2.Rows.Item(4).RowHeight = 18.75

$ws2.Columns.Item("A:M").AutoFit()
$ws2.Columns.Item("A").ColumnWidth = 14.85546875
$ws2.Columns.Item("B").ColumnWidth = 44.85546875
$ws2.Columns.Item("C").ColumnWidth = 17.7109375
$ws2.Columns.Item("D").ColumnWidth = 16.7109375
$ws2.Columns.Item("E:L").ColumnWidth = 6.140625
$ws2.Columns.Item("M").ColumnWidth = 19.42578125

# --- Selection / active sheet restore --------------------------------------

$ws2.Range("O8").Select()
$ws1.Activate()
$ws1.Range("L6").Select()
